$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00213219616204691
$ws.Range("C2").Value = 0.906183368869936
$ws.Range("D2").Value = 0.00426439232409382
$ws.Range("E2").Value = 0.989339019189765
$ws.Range("F2").Value = 0.921108742004264
$ws.Range("H2").Value = 0.00639658848614072
$ws.Range("I2").Value = 0.0383795309168444
$ws.Range("J2").Value = 0.00426439232409382
$ws.Range("K2").Value = 0.889125799573561
$ws.Range("L2").Value = 0.00213219616204691
$ws.Range("N2").Value = 0.00213219616204691
$ws.Range("R2").Value = 0.0149253731343284
$ws.Range("S2").Value = 0.00213219616204691
$ws.Range("T2").Value = 0.997867803837953
$ws.Range("W2").Value = 0.00639658848614072
$ws.Range("B3").Value = 0.991471215351812
$ws.Range("C3").Value = 0.00639658848614072
$ws.Range("D3").Value = 0.991471215351812
$ws.Range("E3").Value = 0.00213219616204691
$ws.Range("F3").Value = 0.00426439232409382
$ws.Range("G3").Value = 0.00213219616204691
$ws.Range("I3").Value = 0.93816631130064
$ws.Range("L3").Value = 0.97228144989339
$ws.Range("M3").Value = 0.023454157782516
$ws.Range("N3").Value = 0.00426439232409382
$ws.Range("O3").Value = 0.997867803837953
$ws.Range("P3").Value = 0.991471215351812
$ws.Range("R3").Value = 0.00852878464818763
$ws.Range("U3").Value = 0.0277185501066098
$ws.Range("V3").Value = 0.00639658848614072
$ws.Range("W3").Value = 0.00213219616204691
$ws.Range("X3").Value = 0.00213219616204691
$ws.Range("B4").Value = 0.00639658848614072
$ws.Range("C4").Value = 0.0874200426439232
$ws.Range("D4").Value = 0.00426439232409382
$ws.Range("E4").Value = 0.00852878464818763
$ws.Range("F4").Value = 0.070362473347548
$ws.Range("G4").Value = 0.995735607675906
$ws.Range("H4").Value = 0.991471215351812
$ws.Range("I4").Value = 0.00852878464818763
$ws.Range("J4").Value = 0.00426439232409382
$ws.Range("K4").Value = 0.00213219616204691
$ws.Range("L4").Value = 0.00639658848614072
$ws.Range("S4").Value = 0.997867803837953
$ws.Range("T4").Value = 0.00213219616204691
$ws.Range("U4").Value = 0.00213219616204691
$ws.Range("V4").Value = 0.00213219616204691
$ws.Range("W4").Value = 0.987206823027719
$ws.Range("X4").Value = 0.993603411513859
$ws.Range("F5").Value = 0.00426439232409382
$ws.Range("G5").Value = 0.00213219616204691
$ws.Range("H5").Value = 0.00213219616204691
$ws.Range("I5").Value = 0.0149253731343284
$ws.Range("J5").Value = 0.991471215351812
$ws.Range("K5").Value = 0.106609808102345
$ws.Range("L5").Value = 0.0191897654584222
$ws.Range("M5").Value = 0.976545842217484
$ws.Range("N5").Value = 0.993603411513859
$ws.Range("P5").Value = 0.00852878464818763
$ws.Range("R5").Value = 0.976545842217484
$ws.Range("U5").Value = 0.970149253731343
$ws.Range("V5").Value = 0.991471215351812
$ws.Range("W5").Value = 0.00426439232409382
$ws.Range("X5").Value = 0.00426439232409382
